$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new row at 100 (shifts old 100-169 down to 101-170)
$ws.Rows.Item(100).Insert()

# Step 2: copy formatting only from row 99 (above, untouched by shift) into row 100
$src = $ws.Range("A99:Q99")
$dst = $ws.Range("A100:Q100")
$src.Copy()
$dst.PasteSpecial(-4122) # xlPasteFormats

$c = $ws.Cells.Item(100, 17) # Q100
Write-Output ("Q100 numfmt=" + $c.NumberFormat)
